$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (single decimal point) need to be
# forced to Text format first, otherwise Excel auto-converts the assigned string into
# a numeric value, which would change the cell type away from the original inline string.
$numericLookingCells = @("D5", "D11", "D15", "D17", "D18", "D25", "D31", "D36", "D39", "D40", "D41", "D42", "D44", "D46", "D48", "D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.383.27'
$ws.Range("E2").Value = '  +4.01%  '
$ws.Range("D3").Value = '1.586.98'
$ws.Range("E3").Value = '  +1.19%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = '214.24'
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("E8").Value = '  +7.58%  '
$ws.Range("E9").Value = '  +0.82%  '
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").Value = '0.0889'
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("D12").Value = '1.813.25'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("D13").Value = '1.590.72'
$ws.Range("E13").Value = '  +1.42%  '
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = '0.532'
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").Value = '28.393.55'
$ws.Range("E16").Value = '  +4.23%  '
$ws.Range("D17").Value = '63.22'
$ws.Range("E17").Value = '  +1.53%  '
$ws.Range("D18").Value = '231.81'
$ws.Range("E18").Value = '  +6.57%  '
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").Value = '152.03'
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").Value = '0.0472'
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("D34").Value = '1.407.10'
$ws.Range("E34").Value = '  -2.67%  '
$ws.Range("E35").Value = '  -1.44%  '
$ws.Range("D36").Value = '1.06'
$ws.Range("E36").Value = '  -4.52%  '
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("D39").Value = '2.53'
$ws.Range("E39").Value = '  +8.34%  '
$ws.Range("D40").Value = '0.542'
$ws.Range("E40").Value = '  +1.66%  '
$ws.Range("D41").Value = '0.817'
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("D42").Value = '5.77'
$ws.Range("E42").Value = '  -2.03%  '
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").Value = '0.982'
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("E45").Value = '  +5.50%  '
$ws.Range("D46").Value = '64.56'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '1.727.09'
$ws.Range("D48").Value = '87.64'
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("E49").Value = '  +5.51%  '
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("D51").Value = '39.22'
$ws.Range("E51").Value = '  +15.69%  '

# Restore the default (General/Normal) style on the cells we temporarily reformatted,
# so the only lasting change is the cell text/value, matching the original formatting.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
